$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the three new header cells, matching the style of the existing
# header row (bold/centered/bordered style) by copying an existing
# header cell's formatting onto the new ones.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record (Wins/Losses/Ties) for every data row (2-43)
$lastRow = 43
$ws.Range("AD2:AD$lastRow").Value = 81
$ws.Range("AE2:AE$lastRow").Value = 81
$ws.Range("AF2:AF$lastRow").Value = 0
